$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A28").Value = "TEST"
$ws.Range("A28").Font.Color = 0
Write-Output "ok"
